$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the source data which stores prices/changes as text strings),
# otherwise Excel auto-converts e.g. "302.58" into the numeric value 302.58.
# We do this by temporarily marking the cell as Text format, assigning the
# value, then resetting the cell style back to Normal so no stray number
# format / style index is left behind on the cell.

$ws.Range("D2").Value = '43.188.92'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '2.323.49'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +2.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.96%  '
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").Value = '2.684.82'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").Value = '2.319.89'
$ws.Range("E16").Value = '  +1.85%  '
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '43.106.33'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.69%  '
$ws.Range("E20").Value = '  +1.51%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("E24").Value = '  +3.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +3.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.65'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.15%  '
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0698'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").Value = '2.003.12'
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("E43").Value = '  +1.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +1.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '74.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.29%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.549.80'
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("E51").Value = '  +2.00%  '
